# Add a new weekly record at the top of the Ciboulette price history block.
# The new row is inserted as row 559, pushing the existing rows 559:669 down
# to 560:670. The inserted row starts as a duplicate of the (old) row 559,
# then the date (Fecha) and volume (Volumen) are updated to the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 559 and insert the copy above itself, shifting the rest of
# the data block (559:669) down by one row (to 560:670).
$ws.Rows.Item(559).Copy()
$ws.Rows.Item(559).Insert()

# Set the new record's own Fecha (D) and Volumen (J) values.
$ws.Cells.Item(559, 4).Value = 45275
$ws.Cells.Item(559, 10).Value = 430
